# The deck currently uses the "Integral" theme (ppt/theme/theme2.xml,
# wired to the slide master / presentation) for its live design.
# The authored change swaps the live theme's 12-color palette back to
# the stock "Office Theme" colors (the palette that used to live,
# unused, in ppt/theme/theme1.xml).
#
# PowerPoint's object model exposes the active theme's 12 theme colors
# through Slide.ThemeColorScheme (indices 1-12, in the standard OOXML
# clrScheme order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink).
# Writing to it edits the underlying theme part in place, which is the
# supported, non-destructive way to recolor the design from COM.

function ConvertTo-OleColor($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Stock "Office Theme" color scheme, in clrScheme slot order.
$officeThemeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme

for ($i = 1; $i -le 12; $i++) {
    $themeColors.Item($i).RGB = ConvertTo-OleColor($officeThemeColors[$i - 1])
}
